# Update "想去人数" (F column) figures on the 展览 and 全部类型 sheets.
# Maps worksheet name -> { row = @(oldValue, newValue) } for documentation/verification.

$wb = $excel.ActiveWorkbook

$updates = @{
    "展览" = @{
        3  = 64
        5  = 14
        6  = 199
        7  = 4626
        14 = 195
        16 = 83
        22 = 3611
        23 = 5949
        37 = 273
        38 = 356
        40 = 1531
        41 = 918
        43 = 33
        47 = 69
    }
    "全部类型" = @{
        3  = 64
        5  = 14
        6  = 199
        7  = 4626
        15 = 195
        17 = 83
        23 = 3611
        24 = 5949
        38 = 273
        39 = 356
        41 = 1531
        42 = 918
        44 = 33
        48 = 69
    }
}

foreach ($sheetName in $updates.Keys) {
    $ws = $wb.Worksheets.Item($sheetName)
    $rowsForSheet = $updates[$sheetName]
    foreach ($row in $rowsForSheet.Keys) {
        $newValue = $rowsForSheet[$row]
        $ws.Cells.Item($row, 6).Value = $newValue
    }
}
